$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current last data row (74), pushing the
# existing row 74 down to row 76 (its contents/formatting move with it).
$ws.Rows.Item(74).Insert()
$ws.Rows.Item(74).Insert()

# Row 74: new "Primera" quality record for 2021-09-09 (serial 44448)
$ws.Range("A74").Value = 1
$ws.Range("B74").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C74").Value = "Arica y Parinacota"
$ws.Range("D74").Value = 44448
$ws.Range("E74").Value = 15
$ws.Range("F74").Value = "Fruta"
$ws.Range("G74").Value = 100108
$ws.Range("H74").Value = "Tropicales y subtropicales"
$ws.Range("I74").Value = 100108003
$ws.Range("J74").Value = "Maracuyá"
$ws.Range("K74").Value = "Sin especificar"
$ws.Range("L74").Value = "Primera"
$ws.Range("M74").Value = 120
$ws.Range("N74").Value = 23000
$ws.Range("O74").Value = 24000
$ws.Range("P74").Value = 23500
$ws.Range("Q74").Value = "$/caja 20 kilos"
$ws.Range("R74").Value = "Región de Arica y Parinacota"
$ws.Range("S74").Value = 1175
$ws.Range("T74").Value = 20

# Row 75: new "Segunda" quality record for 2021-09-09 (serial 44448)
$ws.Range("A75").Value = 1
$ws.Range("B75").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C75").Value = "Arica y Parinacota"
$ws.Range("D75").Value = 44448
$ws.Range("E75").Value = 15
$ws.Range("F75").Value = "Fruta"
$ws.Range("G75").Value = 100108
$ws.Range("H75").Value = "Tropicales y subtropicales"
$ws.Range("I75").Value = 100108003
$ws.Range("J75").Value = "Maracuyá"
$ws.Range("K75").Value = "Sin especificar"
$ws.Range("L75").Value = "Segunda"
$ws.Range("M75").Value = 120
$ws.Range("N75").Value = 20000
$ws.Range("O75").Value = 21000
$ws.Range("P75").Value = 20500
$ws.Range("Q75").Value = "$/caja 20 kilos"
$ws.Range("R75").Value = "Región de Arica y Parinacota"
$ws.Range("S75").Value = 1025
$ws.Range("T75").Value = 20
